# feat: add 2022-Q3 data
#
# Inserts a brand new "2022-Q3" sheet (fund-holdings detail) right after
# the "总计" (totals) sheet, pushing every existing quarter sheet one
# slot to the right, and adds a corresponding new row at the top of the
# "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet immediately after "总计" (i.e.
#    immediately before the sheet that is currently "2022-Q2").
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($anchor)
$q3.Name = "2022-Q3"

# Every quarterly detail sheet shares the same layout/look (bold,
# bordered, centered header row + index column). Clone that formatting
# from the "2022-Q2" sheet instead of rebuilding it property-by-property
# (assigning font/border/alignment piecemeal onto brand-new cells does
# not reliably round-trip in this host).
$src = $wb.Worksheets.Item("2022-Q2")
$src.Range("B1:H1").Copy()
$q3.Range("B1").PasteSpecial(-4122)
$src.Range("A2:H4").Copy()
$q3.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row.
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Force the text-looking columns (fund code + the numeric-looking ratio
# columns D:G) to be stored as text, matching the source data which
# keeps these as strings (e.g. "011486", "9.52") instead of numbers -
# critical for B so the leading zero in fund codes survives.
$q3.Range("B2:B4").NumberFormat = "@"
$q3.Range("D2:G4").NumberFormat = "@"

# Row 2 - 011486 / 博时创新精选混合A
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "011486"
$q3.Range("C2").Value = "博时创新精选混合A"
$q3.Range("D2").Value = "9.52"
$q3.Range("E2").Value = "93.58"
$q3.Range("F2").Value = "2.70"
$q3.Range("G2").Value = "0.2570"
$q3.Range("H2").Value = 10

# Row 3 - 003655 / 信澳新财富灵活配置混合
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "003655"
$q3.Range("C3").Value = "信澳新财富灵活配置混合"
$q3.Range("D3").Value = "4.04"
$q3.Range("E3").Value = "54.87"
$q3.Range("F3").Value = "1.70"
$q3.Range("G3").Value = "0.0687"
$q3.Range("H3").Value = 10

# Row 4 - 011487 / 博时创新精选混合C
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "011487"
$q3.Range("C4").Value = "博时创新精选混合C"
$q3.Range("D4").Value = "1.03"
$q3.Range("E4").Value = "93.58"
$q3.Range("F4").Value = "2.70"
$q3.Range("G4").Value = "0.0278"
$q3.Range("H4").Value = 10

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: a new row for 2022-Q3 is inserted
#    at the top of the data (row 2) and every existing quarter's row
#    shifts down by one. Rewrite the whole data block (rows 2-9) in its
#    final order so the index column (A) and styling stay consistent,
#    rather than relying on a row-insert (which mangles column-A style
#    inheritance).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# The table grows from 8 data rows to 9: extend column A's existing
# "index" formatting (bold/border/centered) one row further down before
# rewriting the values, by cloning the format of the last existing data
# row (row 8) onto the new row (row 9).
$total.Range("A8").Copy()
$total.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rows = @(
    @("2022-Q3", 3, 0.35),
    @("2022-Q2", 6, 1.41),
    @("2022-Q1", 1, 1.56),
    @("2021-Q4", 4, 3.74),
    @("2021-Q3", 3, 0.6),
    @("2021-Q2", 1, 0.58),
    @("2021-Q1", 2, 1.38),
    @("2020-Q4", 3, 2.59)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $data[0]
    $total.Cells.Item($r, 3).Value = $data[1]
    $total.Cells.Item($r, 4).Value = $data[2]
}
